# Update column C ("Förändrad") date values from 2023-10-05 (45204) to
# 2023-10-08 (45207) for rows 2 through 32 on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 32; $row++) {
    $cell = $ws.Cells.Item($row, 3)  # Column C
    if ($cell.Value2 -eq 45204) {
        $cell.Value2 = 45207
    }
}
